# The author added a new weekly price record for "Femacal de La Calera -
# Zanahoria" (Coquimbo). In the sheet's canonical OOXML this shows up as a
# brand-new row 429, with every row that used to be 429..480 shifted down
# by one to become 430..481 (dimension grows from A1:R480 to A1:R481).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 429; Excel automatically shifts
# rows 429..480 down to 430..481 and extends the used range accordingly.
$ws.Rows.Item(429).Insert()

# Populate the newly inserted row 429 with the new record's data.
$ws.Range("A429").Value = 3
$ws.Range("B429").Value = "Femacal de La Calera"
$ws.Range("C429").Value = "Coquimbo"
$ws.Range("D429").Value = 44946
$ws.Range("E429").Value = 5
$ws.Range("F429").Value = 100114013
$ws.Range("G429").Value = "Zanahoria"
$ws.Range("H429").Value = "Sin especificar"
$ws.Range("I429").Value = "Primera"
$ws.Range("J429").Value = 570
$ws.Range("K429").Value = 10000
$ws.Range("L429").Value = 11000
$ws.Range("M429").Value = 10561
$ws.Range("N429").Value = '$/saco 20 kilos'
$ws.Range("O429").Value = "Provincia de Quillota"
$ws.Range("P429").Value = 528
$ws.Range("Q429").Value = 20
$ws.Range("R429").Value = "Hortaliza"
